$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 and F2 values on the remaining row
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.057912333

# Delete rows 3 and 4 entirely (they contained the old puzzle 2 and 3 data)
$ws.Range("A3:F4").Delete()
